# Definition/acronyms/abbreviations document update.
#
# 1. Insert three new paragraphs right before the "Standard customer - ..."
#    definition:
#       "The 4 taxi ride status -"   (with grammar-check marks around "4 taxi")
#       "Taxi driver status -"
#       an otherwise-empty paragraph that now carries the "_GoBack" bookmark
# 2. The hidden "_GoBack" bookmark used to sit in the "Abbreviations"
#    heading paragraph; Word always keeps only one such bookmark, tracking
#    the most recent edit location, so it must be removed from its old spot
#    (it gets recreated at the new empty paragraph above).

$d = $word.ActiveDocument

# --- Step 1: drop the pre-existing "_GoBack" bookmark -----------------
# It is hidden from normal enumeration, but Exists()/Item() still reach it.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 2: locate the "Standard customer" definition paragraph ------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Standard customer*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    Write-Output "ERROR: 'Standard customer' paragraph not found"
} else {
    # Insert three blank paragraphs directly above it, preserving order.
    $target.Range.InsertParagraphBefore()
    $target.Range.InsertParagraphBefore()
    $target.Range.InsertParagraphBefore()

    # Re-locate the (now shifted) target paragraph's index so we can grab
    # the three freshly-created blank paragraphs right before it.
    $targetIdx = -1
    $idx = 0
    foreach ($p in $d.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text -like "Standard customer*") {
            $targetIdx = $idx
            break
        }
    }

    $p1 = $d.Paragraphs.Item($targetIdx - 3)
    $p2 = $d.Paragraphs.Item($targetIdx - 2)
    $p3 = $d.Paragraphs.Item($targetIdx - 1)

    # Give each blank paragraph placeholder text first - InsertXML on a
    # truly empty (collapsed) range does not reliably stick, but replacing
    # a non-empty range with a full <w:p> fragment does.
    $p1.Range.Text = "x"
    $p2.Range.Text = "x"
    $p3.Range.Text = "x"

    $null = $p1.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>4 taxi</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> ride status &#8211;</w:t></w:r></w:p>')

    $null = $p2.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Taxi driver status &#8211;</w:t></w:r></w:p>')

    $null = $p3.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')

    Write-Output "Inserted new definitions and moved the _GoBack bookmark."
}
